$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-08"

# Update the October label in column A (row 11)
$ws.Range("A11").Value = "October (through 10-08)"

# Update September row (row 9, column I = 2022)
$ws.Range("I9").Value = 164

# Update October row already-through values (row 10, column I = 2022)
$ws.Range("I10").Value = 146

# Update October row (row 11) with new carjacking counts
$ws.Range("C11").Value = 15
$ws.Range("D11").Value = 12
$ws.Range("E11").Value = 19
$ws.Range("F11").Value = 8
$ws.Range("H11").Value = 57
$ws.Range("I11").Value = 28

# Update Total row (row 12) to reflect new sums
$ws.Range("C12").Value = 444
$ws.Range("D12").Value = 639
$ws.Range("E12").Value = 567
$ws.Range("F12").Value = 430
$ws.Range("H12").Value = 1304
$ws.Range("I12").Value = 1306
